$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "Which was released first, The Greatest Game Ever Played or Hot Lead and Cold Feet?"
$ws.Range("F2").Value = "Y"
$ws.Range("G2").Value = "Hot Lead and Cold Feet (1977) was released first, preceding The Greatest Game Ever Played (2007) by 30 years. The Western comedy Hot Lead and Cold Feet, starring Roy Rogers, premiered on October 14, 1..."
$ws.Range("H2").Value = "Hot Lead and Cold Feet (1978) was released first, preceding The Greatest Game Ever Played (2005) by 27 years. Hot Lead and Cold Feet, a Western comedy, premiered on July 5, 1978. The Greatest Game Eve..."
$ws.Range("I2").Value = "In evaluating the responses from both assistants, several factors come into play. `nAssistant A correctly identifies that ""Hot Lead and Cold Feet"" was released first and provides the correct release date of October 14, 1977, for that film. However, it incorrectly states that ""The Greatest Game Ever Played"" was released on September 15, 2007, which is inaccurate. The assistant also mentions the different production contexts and genres, which adds some depth to the response.`nAssistant B, on the other hand, incorrectly states the release year of ""Hot Lead and Cold Feet"" as 1978 and the release date of ""The Greatest Game Ever Played"" as September 30, 2005, which is also incorrect. However, it does attempt to clarify the timeline and references verified sources, which adds a layer of credibility to its response.`nIn terms of helpfulness, relevance, and correctness, Assistant A provides a more accurate release date for ""Hot Lead and Cold Feet,"" despite the error regarding ""The Greatest Game Ever Played."" Assistant B's attempt to reference sources is commendable, but the inaccuracies in both release dates diminish its overall effectiveness.`nConsidering these factors, Assistant A's response is more helpful and relevant, despite its own inaccuracies. Therefore, the final verdict is:`n[[A]]"
$ws.Range("J2").Value = "A"
$ws.Range("P2").Value = "Y"

# --- Row 3 ---
$ws.Range("B3").Value = "Which was released first, The Greatest Game Ever Played or Hot Lead and Cold Feet?"
$ws.Range("G3").Value = """Hot Lead and Cold Feet"" was released first. It is a Disney Western comedy film released in 1978. ""The Greatest Game Ever Played"" is a golf drama that was released much later, in 2005...."
$ws.Range("H3").Value = "Hot Lead and Cold Feet was released first, on July 5, 1978. The Greatest Game Ever Played was released on September 30, 2005. `nReferences:`n- https://en.wikipedia.org/wiki/Hot_Lead_and_Cold_Feet`n- htt..."
$ws.Range("I3").Value = "Both Assistant A and Assistant B provide the correct information regarding the release dates of ""Hot Lead and Cold Feet"" and ""The Greatest Game Ever Played."" However, there are notable differences in their responses.`nAssistant A states that ""Hot Lead and Cold Feet"" was released first and provides the year of release for both films. While this answer is correct, it lacks specific details such as the exact release date of ""Hot Lead and Cold Feet,"" which could enhance the depth of the response.`nAssistant B also correctly identifies that ""Hot Lead and Cold Feet"" was released first and includes the exact release date of July 5, 1978, as well as the release date of ""The Greatest Game Ever Played"" on September 30, 2005. Additionally, Assistant B provides references to credible sources, which adds to the accuracy and reliability of the information presented.`nIn terms of helpfulness, relevance, correctness, and depth, Assistant B outperforms Assistant A by providing more detailed information and supporting references.`nFinal verdict: [[B]]"
